$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text type (so numeric-looking strings like "1.000" or "29.082.21"
# are not reinterpreted as numbers/dates), then strip the resulting
# number-format override so cell styling matches the original (unstyled) cells.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.082.21'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '1.820.70'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("D4").Value = '0.9982'
$ws.Range("D5").Value = '241.47'
$ws.Range("E5").Value = '  -0.77%  '
$ws.Range("E6").Value = '  -2.16%  '
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  -2.46%  '
$ws.Range("D9").Value = '0.2875'
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("E10").Value = '  -1.48%  '
$ws.Range("D11").Value = '0.07655'
$ws.Range("E11").Value = '  -0.33%  '
$ws.Range("D12").Value = '1.822.52'
$ws.Range("E12").Value = '  -0.63%  '
$ws.Range("E13").Value = '  -1.35%  '
$ws.Range("D14").Value = '0.6591'
$ws.Range("E14").Value = '  -1.12%  '
$ws.Range("D15").Value = '81.42'
$ws.Range("E15").Value = '  -1.53%  '
$ws.Range("D16").Value = '0.000008956'
$ws.Range("E16").Value = '  -4.42%  '
$ws.Range("D17").Value = '5.826'
$ws.Range("E17").Value = '  -2.49%  '
$ws.Range("D18").Value = '29.061.65'
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").Value = '2.058.31'
$ws.Range("E19").Value = '  -1.12%  '
$ws.Range("D20").Value = '237.89'
$ws.Range("E20").Value = '  +6.58%  '
$ws.Range("E21").Value = '  -1.30%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("D23").Value = '7.104'
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").Value = '0.9997'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").Value = '157.68'
$ws.Range("D26").Value = '0.1402'
$ws.Range("E26").Value = '  +0.70%  '
$ws.Range("D27").Value = '8.405'
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("D28").Value = '17.55'
$ws.Range("E28").Value = '  -1.99%  '
$ws.Range("D29").Value = '1.477'
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("D30").Value = '0.05561'
$ws.Range("E30").Value = '  -2.17%  '
$ws.Range("D31").Value = '4.084'
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("D32").Value = '4.092'
$ws.Range("E32").Value = '  -1.35%  '
$ws.Range("E33").Value = '  +0.21%  '
$ws.Range("D34").Value = '0.7335'
$ws.Range("E34").Value = '  -1.21%  '
$ws.Range("D35").Value = '1.810'
$ws.Range("E35").Value = '  -1.53%  '
$ws.Range("E36").Value = '  -1.09%  '
$ws.Range("D37").Value = '2.620'
$ws.Range("E37").Value = '  -1.88%  '
$ws.Range("D38").Value = '2.826'
$ws.Range("E38").Value = '  +2.31%  '
$ws.Range("D39").Value = '1.206.07'
$ws.Range("E39").Value = '  -1.28%  '
$ws.Range("D40").Value = '0.01756'
$ws.Range("E40").Value = '  -1.30%  '
$ws.Range("D41").Value = '6.347'
$ws.Range("E41").Value = '  -2.88%  '
$ws.Range("D42").Value = '0.8907'
$ws.Range("E42").Value = '  -0.35%  '
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").Value = '100.75'
$ws.Range("E44").Value = '  -1.14%  '
$ws.Range("D45").Value = '1.965.65'
$ws.Range("E45").Value = '  -0.83%  '
$ws.Range("D46").Value = '64.47'
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.00000000122'
$ws.Range("E47").Value = '  -2.62%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '0.5081'
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("D49").Value = '9.040'
$ws.Range("E49").Value = '  +0.68%  '
$ws.Range("E50").Value = '  -2.39%  '
$ws.Range("D51").Value = '0.05749'
$ws.Range("E51").Value = '  -1.12%  '

$editRange.ClearFormats()
